# Insert a new weekly price record at row 373 for
# "Terminal La Palmera de La Serena - Poroto verde".
# This pushes the former rows 373-422 down to 374-423
# (Excel's native row-insert behavior handles the shift,
# formatting/style inheritance, and dimension update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(373).Insert()

$ws.Cells.Item(373, 1).Value  = 8
$ws.Cells.Item(373, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(373, 3).Value  = "Coquimbo"
$ws.Cells.Item(373, 4).Value  = 45131
$ws.Cells.Item(373, 5).Value  = 4
$ws.Cells.Item(373, 6).Value  = 100112031
$ws.Cells.Item(373, 7).Value  = "Poroto verde"
$ws.Cells.Item(373, 8).Value  = "Magnum"
$ws.Cells.Item(373, 9).Value  = "Primera"
$ws.Cells.Item(373, 10).Value = 460
$ws.Cells.Item(373, 11).Value = 24000
$ws.Cells.Item(373, 12).Value = 25000
$ws.Cells.Item(373, 13).Value = 24500
$ws.Cells.Item(373, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(373, 15).Value = "Perú"
$ws.Cells.Item(373, 16).Value = 980
$ws.Cells.Item(373, 17).Value = 25
$ws.Cells.Item(373, 18).Value = "Hortaliza"
